$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 20, shifting the existing rows 20-99 down to 21-100.
$ws.Rows.Item(20).Insert()

# Copy the date-format style from the row above (row 19's D column) onto the
# new row's D cell so it keeps the same date display.
$ws.Range("D19").Copy()
$ws.Range("D20").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new record's data (same static columns as the rest of the block,
# new values for the date/volume/price columns).
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = "Vega Modelo de Temuco"
$ws.Range("C20").Value = "La Araucanía"
$ws.Range("D20").Value = 45069
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100108
$ws.Range("H20").Value = "Tropicales y subtropicales"
$ws.Range("I20").Value = 100108007
$ws.Range("J20").Value = "Coco"
$ws.Range("K20").Value = "Sin especificar"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 8
$ws.Range("N20").Value = 36000
$ws.Range("O20").Value = 36000
$ws.Range("P20").Value = 36000
$ws.Range("Q20").Value = "$/malla 20 unidades"
$ws.Range("R20").Value = "Perú"
$ws.Range("S20").Value = 1800
$ws.Range("T20").Value = 20
